$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 6204.8
$ws.Range("J2").Value = 10002
$ws.Range("L2").Value = 10002
$ws.Range("N2").Value = -10228
$ws.Range("H28").Value = 104.5
$ws.Range("I28").Value = 108
$ws.Range("K28").Value = 108
$ws.Range("M28").Value = 377
$ws.Range("H52").Value = 396.66666
$ws.Range("I52").Value = 396.66666
$ws.Range("K52").Value = 1189.99998
$ws.Range("M52").Value = -1029.99998
$ws.Range("H55").Value = 1198.5
$ws.Range("I55").Value = 1041.5
$ws.Range("J55").Value = 1394.75
$ws.Range("K55").Value = 1041.5
$ws.Range("L55").Value = 1394.75
$ws.Range("M55").Value = -827.5
$ws.Range("N55").Value = -1822.75
$ws.Range("H62").Value = 220004700
$ws.Range("I62").Value = 275003360
$ws.Range("K62").Value = 275003360
$ws.Range("M62").Value = -275002736
$ws.Range("H64").Value = 5042.375
$ws.Range("I64").Value = 4605.143
$ws.Range("K64").Value = 4605.143
$ws.Range("M64").Value = -4357.143
$ws.Range("H65").Value = 220004700
$ws.Range("I65").Value = 275003360
$ws.Range("K65").Value = 1375016800
$ws.Range("M65").Value = -1375013680
$ws.Range("H67").Value = 5042.375
$ws.Range("I67").Value = 4605.143
$ws.Range("K67").Value = 4605.143
$ws.Range("M67").Value = -3747.143
$ws.Range("H105").Value = 44999.5
$ws.Range("J105").Value = 44999.5
$ws.Range("L105").Value = 44999.5
$ws.Range("N105").Value = -51987.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9433.700000000001
$ws.Range("I61").Value = 4258.9165
$ws.Range("J61").Value = 17195.875
$ws.Range("K61").Value = 4258.9165
$ws.Range("L61").Value = 17195.875
$ws.Range("M61").Value = -4046.9165
$ws.Range("N61").Value = -17619.875
$ws.Range("H76").Value = 20288
$ws.Range("J76").Value = 20288
$ws.Range("L76").Value = 20288
$ws.Range("N76").Value = -20964
$ws.Range("H79").Value = 20288
$ws.Range("J79").Value = 20288
$ws.Range("L79").Value = 20288
$ws.Range("N79").Value = -22628
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26802
$ws.Range("H103").Value = 39787.332
$ws.Range("J103").Value = 39787.332
$ws.Range("L103").Value = 39787.332
$ws.Range("N103").Value = -42131.332
$ws.Range("H107").Value = 86074.336
$ws.Range("J107").Value = 86074.336
$ws.Range("L107").Value = 86074.336
$ws.Range("N107").Value = -93754.336
$ws.Range("H136").Value = 9433.700000000001
$ws.Range("I136").Value = 4258.9165
$ws.Range("J136").Value = 17195.875
$ws.Range("K136").Value = 12776.7495
$ws.Range("L136").Value = 51587.625
$ws.Range("M136").Value = -10226.7495
$ws.Range("N136").Value = -56687.625

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 24098.111
$ws.Range("I54").Value = 3376.6
$ws.Range("K54").Value = 3376.6
$ws.Range("M54").Value = -2892.6
$ws.Range("H108").Value = 122676.5
$ws.Range("J108").Value = 122676.5
$ws.Range("L108").Value = 122676.5
$ws.Range("N108").Value = -130356.5
$ws.Range("H111").Value = 98494.5
$ws.Range("J111").Value = 98494.5
$ws.Range("L111").Value = 98494.5
$ws.Range("N111").Value = -106674.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 689.1053000000001
$ws.Range("I22").Value = 593.7646999999999
$ws.Range("K22").Value = 593.7646999999999
$ws.Range("M22").Value = -243.7646999999999
$ws.Range("H92").Value = 58901
$ws.Range("J92").Value = 58901
$ws.Range("L92").Value = 58901
$ws.Range("N92").Value = -63893

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 5575.25
$ws.Range("I9").Value = 7100.3335
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 21301.0005
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = -21077.0005
$ws.Range("N9").Value = -3448
$ws.Range("H12").Value = 109.26667
$ws.Range("J12").Value = 124.23077
$ws.Range("L12").Value = 372.69231
$ws.Range("N12").Value = -718.69231
$ws.Range("H26").Value = 3212.2144
$ws.Range("I26").Value = 2354
$ws.Range("J26").Value = 4757
$ws.Range("K26").Value = 7062
$ws.Range("L26").Value = 14271
$ws.Range("M26").Value = -6774
$ws.Range("N26").Value = -14847
$ws.Range("H68").Value = 1861.6364
$ws.Range("J68").Value = 2292.875
$ws.Range("L68").Value = 6878.625
$ws.Range("N68").Value = -8500.625
$ws.Range("H71").Value = 1861.6364
$ws.Range("J71").Value = 2292.875
$ws.Range("L71").Value = 20635.875
$ws.Range("N71").Value = -28747.875
$ws.Range("H75").Value = 599.6
$ws.Range("J75").Value = 749
$ws.Range("L75").Value = 2247
$ws.Range("N75").Value = -4243
$ws.Range("H78").Value = 599.6
$ws.Range("J78").Value = 749
$ws.Range("L78").Value = 6741
$ws.Range("N78").Value = -16725
$ws.Range("H107").Value = 750.0952
$ws.Range("I107").Value = 414.23077
$ws.Range("J107").Value = 1295.875
$ws.Range("K107").Value = 1242.69231
$ws.Range("L107").Value = 3887.625
$ws.Range("M107").Value = 677.3076900000001
$ws.Range("N107").Value = -7727.625
$ws.Range("H140").Value = 19232510
$ws.Range("J140").Value = 2236.8572
$ws.Range("L140").Value = 6710.571599999999
$ws.Range("N140").Value = -17070.5716

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 26000000
$ws.Range("I3").Value = 26000000
$ws.Range("K3").Value = 26000000
$ws.Range("M3").Value = -25999884
$ws.Range("H43").Value = 6596.8
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20302
$ws.Range("H46").Value = 49157.75
$ws.Range("J46").Value = 55489.3
$ws.Range("L46").Value = 55489.3
$ws.Range("N46").Value = -55801.3
$ws.Range("H92").Value = 5000
$ws.Range("J92").Value = 5000
$ws.Range("L92").Value = 5000
$ws.Range("N92").Value = -8744
$ws.Range("H106").Value = 60000
$ws.Range("J106").Value = 60000
$ws.Range("L106").Value = 60000
$ws.Range("N106").Value = -62524
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 500
$ws.Range("J122").Value = 4250
$ws.Range("K122").Value = 1500
$ws.Range("L122").Value = 12750
$ws.Range("M122").Value = 950
$ws.Range("N122").Value = -17650
$ws.Range("H134").Value = 78328.81
$ws.Range("J134").Value = 78328.81
$ws.Range("L134").Value = 234986.43
$ws.Range("N134").Value = -240056.43

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1587.2858
$ws.Range("I82").Value = 1750
$ws.Range("J82").Value = 1370.3334
$ws.Range("K82").Value = 1750
$ws.Range("L82").Value = 1370.3334
$ws.Range("M82").Value = -1389
$ws.Range("N82").Value = -2092.3334
$ws.Range("H85").Value = 1587.2858
$ws.Range("I85").Value = 1750
$ws.Range("J85").Value = 1370.3334
$ws.Range("K85").Value = 1750
$ws.Range("L85").Value = 1370.3334
$ws.Range("M85").Value = -502
$ws.Range("N85").Value = -3866.3334
$ws.Range("H93").Value = 8003
$ws.Range("I93").Value = 8003
$ws.Range("K93").Value = 8003
$ws.Range("M93").Value = -6755
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""   # was -31352
$ws.Range("H103").Value = 22896.25
$ws.Range("J103").Value = 22896.25
$ws.Range("L103").Value = 22896.25
$ws.Range("N103").Value = -25240.25
$ws.Range("H122").Value = 21743800
$ws.Range("I122").Value = 27782546
$ws.Range("K122").Value = 83347638
$ws.Range("M122").Value = -83345188
$ws.Range("H135").Value = 141666.5
$ws.Range("J135").Value = 141666.5
$ws.Range("L135").Value = 141666.5
$ws.Range("N135").Value = -151806.5
$ws.Range("H136").Value = 8922.85
$ws.Range("I136").Value = 7850.4116
$ws.Range("K136").Value = 23551.2348
$ws.Range("M136").Value = -21001.2348

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""   # was -57680
$ws.Range("H141").Value = 119994
$ws.Range("J141").Value = 119994
$ws.Range("L141").Value = 119994
$ws.Range("N141").Value = -130354
